$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 469
$ws1.Range("F3").Value = 5621
$ws1.Range("F6").Value = 87
$ws1.Range("F9").Value = 532

# Sheet "全部类型" (All types) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 469
$ws4.Range("F3").Value = 5621
$ws4.Range("F7").Value = 87
$ws4.Range("F11").Value = 532
